# Apply text corrections described in the diff (v1.3 -> v1.4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) TC2 step: "excluido" -> "nao excluido" (shared string used by D20)
$ws.Range("D20").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos com o Periodo Avaliativo nao excluido"

# 2) Wording fix for "Data Inicial e Data Final" -> "Data Inicial' e 'Data Final"
# These strings are reused across several rows (TC5-TC10 step rows)
$rows = @(51, 63, 75, 87, 98, 109)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "Lider de Pessoas preenche o campo 'Data Inicial' e 'Data Final' informando as respectivas datas referentes ao periodo"
    $ws.Range("D$r").Value = "SYSTEM apresenta o campo 'Data Inicial' e 'Data Final' preenchido corretamente"
}
